$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.060.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "'2.886.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.79%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'352.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'111.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").Value = "'0.561"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'40.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "'3.335.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.909.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.992"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.19%  "
$ws.Range("D18").Value = "'52.053.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'7.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "'3.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.78%  "
$ws.Range("D21").Value = "'13.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "'71.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "'270.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "'2.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "'26.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("D30").Value = "'38.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.91%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'6.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "'6.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.27%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'53.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("D35").Value = "'0.0932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.54%  "
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.41%  "
$ws.Range("D39").Value = "'18.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("D43").Value = "'22.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("D44").Value = "'122.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "'3.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.47%  "
$ws.Range("D47").Value = "'2.191.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").Value = "'2.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.12%  "
$ws.Range("D49").Value = "'0.267"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +18.85%  "
$ws.Range("D50").Value = "'0.951"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.88%  "
$ws.Range("D51").Value = "'5.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
